$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $ref, $text)
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-CellText $ws "D2" "41.700.99"
Set-CellText $ws "E2" "  +0.37%  "
Set-CellText $ws "D3" "2.466.36"
Set-CellText $ws "E3" "  -0.89%  "
Set-CellText $ws "E4" "  +0.43%  "
Set-CellText $ws "D5" "315.10"
Set-CellText $ws "E5" "  +0.55%  "
Set-CellText $ws "D6" "92.44"
Set-CellText $ws "E6" "  -1.38%  "
Set-CellText $ws "E7" "  +0.43%  "
Set-CellText $ws "E8" "  +0.36%  "
Set-CellText $ws "E9" "  +2.72%  "
Set-CellText $ws "D10" "32.49"
Set-CellText $ws "E10" "  -1.20%  "
Set-CellText $ws "D11" "0.0844"
Set-CellText $ws "E11" "  +7.47%  "
Set-CellText $ws "D13" "2.844.03"
Set-CellText $ws "E13" "  -1.13%  "
Set-CellText $ws "D14" "6.86"
Set-CellText $ws "E14" "  -0.21%  "
Set-CellText $ws "D15" "15.77"
Set-CellText $ws "E15" "  +1.80%  "
Set-CellText $ws "D16" "2.467.13"
Set-CellText $ws "E16" "  -1.74%  "
Set-CellText $ws "E17" "  +2.20%  "
Set-CellText $ws "D18" "41.664.84"
Set-CellText $ws "E18" "  -0.09%  "
Set-CellText $ws "E19" "  +2.35%  "
Set-CellText $ws "D20" "0.0₃0948"
Set-CellText $ws "E20" "  +2.55%  "
Set-CellText $ws "D21" "70.75"
Set-CellText $ws "E21" "  +0.21%  "
Set-CellText $ws "D22" "11.39"
Set-CellText $ws "E22" "  +1.53%  "
Set-CellText $ws "D23" "238.51"
Set-CellText $ws "E23" "  +0.84%  "
Set-CellText $ws "E25" "  +0.64%  "
Set-CellText $ws "E26" "  +0.15%  "
Set-CellText $ws "D27" "24.46"
Set-CellText $ws "E27" "  -0.78%  "
Set-CellText $ws "D28" "2.26"
Set-CellText $ws "E28" "  +0.55%  "
Set-CellText $ws "D29" "9.74"
Set-CellText $ws "E29" "  +0.63%  "
Set-CellText $ws "D30" "35.21"
Set-CellText $ws "E30" "  -3.24%  "
Set-CellText $ws "D31" "155.79"
Set-CellText $ws "E31" "  +0.97%  "
Set-CellText $ws "D32" "5.50"
Set-CellText $ws "E32" "  +1.55%  "
Set-CellText $ws "E33" "  +0.52%  "
Set-CellText $ws "D34" "0.0760"
Set-CellText $ws "E34" "  +0.55%  "
Set-CellText $ws "E35" "  -0.51%  "
Set-CellText $ws "D36" "17.43"
Set-CellText $ws "E36" "  -5.46%  "
Set-CellText $ws "E37" "  -2.71%  "
Set-CellText $ws "D38" "0.115"
Set-CellText $ws "E38" "  +1.02%  "
Set-CellText $ws "E39" "  +0.48%  "
Set-CellText $ws "E40" "  -2.28%  "
Set-CellText $ws "E41" "  -5.29%  "
Set-CellText $ws "E42" "  +0.47%  "
Set-CellText $ws "D43" "1.973.34"
Set-CellText $ws "E43" "  +1.12%  "
Set-CellText $ws "B44" "EnergySwap"
Set-CellText $ws "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D44" "18.94"
Set-CellText $ws "E44" "  -3.84%  "
Set-CellText $ws "B45" "VeChain"
Set-CellText $ws "C45" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D45" "0.0282"
Set-CellText $ws "E45" "  -1.09%  "
Set-CellText $ws "E46" "  -1.39%  "
Set-CellText $ws "D47" "9.01"
Set-CellText $ws "E47" "  +2.18%  "
Set-CellText $ws "D48" "2.700.42"
Set-CellText $ws "E48" "  -1.31%  "
Set-CellText $ws "D49" "96.76"
Set-CellText $ws "E49" "  +0.33%  "
Set-CellText $ws "D50" "66.99"
Set-CellText $ws "E50" "  -0.97%  "
Set-CellText $ws "B51" "MultiversX"
Set-CellText $ws "C51" "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-CellText $ws "D51" "52.36"
Set-CellText $ws "E51" "  +3.18%  "

Write-Host "Applied 86 cell updates"
